$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 166.7619
$ws.Range("I33").Value = 167.75
$ws.Range("J33").Value = 147
$ws.Range("K33").Value = 167.75
$ws.Range("L33").Value = 147
$ws.Range("M33").Value = 61.25
$ws.Range("N33").Value = -605

$ws.Range("H76").Value = 3596.889
$ws.Range("I76").Value = 3037.1667
$ws.Range("K76").Value = 3037.1667
$ws.Range("M76").Value = -2722.1667

$ws.Range("H79").Value = 3596.889
$ws.Range("I79").Value = 3037.1667
$ws.Range("K79").Value = 3037.1667
$ws.Range("M79").Value = -1945.1667

$ws.Range("H112").Value = 2107.9565
$ws.Range("J112").Value = 2135.5908
$ws.Range("L112").Value = 6406.7724
$ws.Range("N112").Value = -8622.7724

$ws.Range("H127").Value = 549.86664
$ws.Range("I127").Value = 446.35715
$ws.Range("K127").Value = 1339.07145
$ws.Range("M127").Value = 3620.92855

$ws.Range("H132").Value = 9796.429
$ws.Range("I132").Value = 7415.85
$ws.Range("J132").Value = 15747.875
$ws.Range("K132").Value = 22247.55
$ws.Range("L132").Value = 47243.625
$ws.Range("M132").Value = -19717.55
$ws.Range("N132").Value = -52303.625

$ws.Range("H138").Value = 3147.5483
$ws.Range("I138").Value = 4659.2666
$ws.Range("K138").Value = 13977.7998
$ws.Range("M138").Value = -8837.799800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2577.889
$ws.Range("I61").Value = 2072.7
$ws.Range("J61").Value = 3209.375
$ws.Range("K61").Value = 2072.7
$ws.Range("L61").Value = 3209.375
$ws.Range("M61").Value = -1860.7
$ws.Range("N61").Value = -3633.375

$ws.Range("H63").Value = 2683.1667
$ws.Range("J63").Value = 2422.5
$ws.Range("L63").Value = 2422.5
$ws.Range("N63").Value = -3794.5

$ws.Range("H66").Value = 2683.1667
$ws.Range("J66").Value = 2422.5
$ws.Range("L66").Value = 12112.5
$ws.Range("N66").Value = -18976.5

$ws.Range("H132").Value = 3743.1875
$ws.Range("I132").Value = 2488.4443
$ws.Range("J132").Value = 5356.4287
$ws.Range("K132").Value = 7465.3329
$ws.Range("L132").Value = 16069.2861
$ws.Range("M132").Value = -4935.3329
$ws.Range("N132").Value = -21129.2861

$ws.Range("H136").Value = 2577.889
$ws.Range("I136").Value = 2072.7
$ws.Range("J136").Value = 3209.375
$ws.Range("K136").Value = 6218.099999999999
$ws.Range("L136").Value = 9628.125
$ws.Range("M136").Value = -3668.099999999999
$ws.Range("N136").Value = -14728.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2504.4666
$ws.Range("I86").Value = 2248.3333
$ws.Range("J86").Value = 3529
$ws.Range("K86").Value = 2248.3333
$ws.Range("L86").Value = 3529
$ws.Range("M86").Value = -1125.3333
$ws.Range("N86").Value = -5775

$ws.Range("H89").Value = 2504.4666
$ws.Range("I89").Value = 2248.3333
$ws.Range("J89").Value = 3529
$ws.Range("K89").Value = 11241.6665
$ws.Range("L89").Value = 17645
$ws.Range("M89").Value = -5625.666499999999
$ws.Range("N89").Value = -28877

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("N108").ClearContents()

$ws.Range("H134").Value = 2524.1
$ws.Range("I134").Value = 2202.0686
$ws.Range("J134").Value = 3906.9412
$ws.Range("K134").Value = 6606.2058
$ws.Range("L134").Value = 11720.8236
$ws.Range("M134").Value = -4071.2058
$ws.Range("N134").Value = -16790.8236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2889.5
$ws.Range("I58").Value = 1885.3334
$ws.Range("J58").Value = 3421.1177
$ws.Range("K58").Value = 1885.3334
$ws.Range("L58").Value = 3421.1177
$ws.Range("M58").Value = -1682.3334
$ws.Range("N58").Value = -3827.1177

$ws.Range("H99").Value = 5717.1113
$ws.Range("I99").Value = 4337.3335
$ws.Range("K99").Value = 4337.3335
$ws.Range("M99").Value = -2839.3335

$ws.Range("H126").Value = 5717.1113
$ws.Range("I126").Value = 4337.3335
$ws.Range("K126").Value = 13012.0005
$ws.Range("M126").Value = -10542.0005

$ws.Range("H136").Value = 2889.5
$ws.Range("I136").Value = 1885.3334
$ws.Range("J136").Value = 3421.1177
$ws.Range("K136").Value = 5656.0002
$ws.Range("L136").Value = 10263.3531
$ws.Range("M136").Value = -3106.0002
$ws.Range("N136").Value = -15363.3531

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3925033.2
$ws.Range("I4").Value = 3537161.8
$ws.Range("J4").Value = 7803749.5
$ws.Range("K4").Value = 10611485.4
$ws.Range("L4").Value = 23411248.5
$ws.Range("M4").Value = -10611373.4
$ws.Range("N4").Value = -23411472.5

$ws.Range("H60").Value = 1179935.1
$ws.Range("I60").Value = 4000503.5
$ws.Range("K60").Value = 12001510.5
$ws.Range("M60").Value = -12001259.5

$ws.Range("H76").Value = 7499
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 7499
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 22497
$ws.Range("N76").Value = -23263
$ws.Range("M76").ClearContents()

$ws.Range("H79").Value = 7499
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 7499
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 22497
$ws.Range("N79").Value = -25149
$ws.Range("M79").ClearContents()

$ws.Range("H113").Value = 1891.3572
$ws.Range("J113").Value = 1998.3846
$ws.Range("L113").Value = 5995.1538
$ws.Range("N113").Value = -10335.1538

$ws.Range("H131").Value = 2196.8462
$ws.Range("J131").Value = 2074.276
$ws.Range("L131").Value = 6222.828
$ws.Range("N131").Value = -16302.828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2944.348
$ws.Range("I132").Value = 2499.5334
$ws.Range("J132").Value = 3778.375
$ws.Range("K132").Value = 7498.600199999999
$ws.Range("L132").Value = 11335.125
$ws.Range("M132").Value = -4968.600199999999
$ws.Range("N132").Value = -16395.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4590.846
$ws.Range("I136").Value = 4334.727
$ws.Range("K136").Value = 13004.181
$ws.Range("M136").Value = -10454.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6945.875
$ws.Range("I62").Value = 3884
$ws.Range("J62").Value = 7966.5
$ws.Range("K62").Value = 3884
$ws.Range("L62").Value = 7966.5
$ws.Range("M62").Value = -3260
$ws.Range("N62").Value = -9214.5

$ws.Range("H65").Value = 6945.875
$ws.Range("I65").Value = 3884
$ws.Range("J65").Value = 7966.5
$ws.Range("K65").Value = 19420
$ws.Range("L65").Value = 39832.5
$ws.Range("M65").Value = -16300
$ws.Range("N65").Value = -46072.5

$ws.Range("H132").Value = 3843.0908
$ws.Range("I132").Value = 4190.3125
$ws.Range("J132").Value = 2917.1667
$ws.Range("K132").Value = 12570.9375
$ws.Range("L132").Value = 8751.500100000001
$ws.Range("M132").Value = -10040.9375
$ws.Range("N132").Value = -13811.5001
